$wb = $excel.ActiveWorkbook

# --- Sheet "3": fix K6:K29 CONCAT formulas (remove trailing space inside
#     "--subst_type without " literal) and update the view/selection ---
$ws3 = $wb.Worksheets.Item("3")
for ($r = 6; $r -le 29; $r++) {
    $cell = $ws3.Range("K" + $r)
    $cell.Formula = '=_xlfn.CONCAT("python 2_graph_processing/make_graph_data.py", " --subst_type without",  " -dir ", dirs!$A$4, "/", ''3''!A' + $r + ')'
}

# --- Sheet "4": fix F2:F19 CONCAT formulas (drop the redundant " " piece
#     right before " -o ") and update the selection ---
$ws4 = $wb.Worksheets.Item("4")
for ($r = 2; $r -le 19; $r++) {
    $cell = $ws4.Range("F" + $r)
    $cell.Formula = '=_xlfn.CONCAT("python 2_graph_processing/make_main_data_combined.py", " -i ", dirs!$A$4, "/", B' + $r + ', " ", dirs!$A$4, "/", C' + $r + ', " -o ", dirs!$A$4, "/", D' + $r + ', " --subst_type ", E' + $r + ')'
}

$excel.CalculateFull()

# --- View state: move the active tab from sheet "2" to sheet "3" ---
$ws3.Activate()
$ws3.Range("K6:K29").Select()

$ws4.Range("F2:F19").Select()

# Restore sheet "3" as the one left active/selected (matches tabSelected
# moving to sheet "3" / activeTab=3 in the saved workbook).
$ws3.Activate()
$ws3.Range("K6:K29").Select()

$wb.Save()
